# Fixed issue #13 Permitir que en los ficheros de metadatos dos columnas
# se puedan relacionar para crear SKOS jerárquicos
#
# Inserts a new row below the header row containing the "slug" (lower-case,
# hyphenated) identifiers for each column, shifting the previously existing
# rows (measure URIs, type, datatype) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing existing row 2 (and below) down.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the lower-case slug values that
# correspond to the header labels in row 1.
$ws.Range("A2").Value = "edad"
$ws.Range("B2").Value = "personas-residentes-viviendas-familiares"
$ws.Range("C2").Value = "situacion-profesional-codigo"
$ws.Range("D2").Value = "aragon"
$ws.Range("E2").Value = "situacion-profesional"
$ws.Range("F2").Value = "sexo"
